$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32: phone number "09876543" was stored as text (to preserve the
# leading zero). Convert it to a plain number (losing the leading zero) -
# it keeps its existing points value (120).
$ws.Range("A32").Value = 9876543

# New row 33: the same customer re-added with the original text phone
# number (leading zero preserved) and their points reset to 0.
# A leading apostrophe forces text storage (keeping the leading zero)
# instead of Excel auto-coercing the digit string to a number; reset the
# style afterwards so no stray quote-prefix formatting sticks around.
$ws.Range("A33").Value = "'09876543"
$ws.Range("A33").Style = "Normal"

# B33 mirrors B32: an explicit empty-text cell (not merely a blank one).
# Assigning "'" makes Excel store an empty quote-prefixed text value; reset
# the style afterwards too.
$ws.Range("B33").Value = "'"
$ws.Range("B33").Style = "Normal"

$ws.Range("C33").Value = 0
